# Applies the "Updated cryptos list" data refresh to Sheet1 of the workbook.
# For each changed row, the Price (D) and/or Volume(1h) (E) columns are updated
# to the newly scraped values; rows 17/18 additionally swap their Coin/Link
# (B/C) contents because WrappedEther and Polkadot traded ranking positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values that are not at risk of being auto-interpreted as numbers by Excel
# (percentages keep their padding spaces, names/links are non-numeric, and a
# couple of "price" values contain multiple thousand-separator dots). ---
$ws.Range("D2").Value = "56.142.92"
$ws.Range("E2").Value = "  -3.49%  "
$ws.Range("D3").Value = "2.929.53"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  -7.18%  "
$ws.Range("E6").Value = "  -7.64%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -6.21%  "
$ws.Range("E9").Value = "  -7.06%  "
$ws.Range("E10").Value = "  -7.75%  "
$ws.Range("E11").Value = "  -6.30%  "
$ws.Range("D12").Value = "3.442.57"
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("E14").Value = "  -5.87%  "
$ws.Range("E15").Value = "  -10.40%  "
$ws.Range("D16").Value = "56.303.84"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.932.39"
$ws.Range("E17").Value = "  -4.66%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E18").Value = "  -5.16%  "
$ws.Range("E19").Value = "  -6.10%  "
$ws.Range("E20").Value = "  -6.75%  "
$ws.Range("E21").Value = "  -8.34%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  -5.77%  "
$ws.Range("E25").Value = "  -4.81%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("D28").Value = "0.0₃0837"
$ws.Range("E28").Value = "  -13.80%  "
$ws.Range("E29").Value = "  -9.21%  "
$ws.Range("E30").Value = "  -7.82%  "
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").Value = "  -6.59%  "
$ws.Range("E33").Value = "  -9.84%  "
$ws.Range("E34").Value = "  -4.90%  "
$ws.Range("E35").Value = "  -8.69%  "
$ws.Range("E36").Value = "  -6.49%  "
$ws.Range("E37").Value = "  -10.44%  "
$ws.Range("E38").Value = "  -10.32%  "
$ws.Range("E39").Value = "  -7.85%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").Value = "2.962.35"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -8.22%  "
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").Value = "2.121.12"
$ws.Range("E45").Value = "  -9.31%  "
$ws.Range("E46").Value = "  -10.59%  "
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("E48").Value = "  -13.75%  "
$ws.Range("E49").Value = "  -6.75%  "
$ws.Range("E50").Value = "  -7.60%  "
$ws.Range("E51").Value = "  -7.11%  "

# --- Price values that look like plain decimal numbers. These cells store text
# (e.g. "0.420", "1.00") in the workbook, so we briefly force a Text number
# format before assigning, then restore General, to stop Excel from silently
# converting them to numbers and dropping significant trailing zeros. ---
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "490.75"
$ws.Range("D6").Value = "132.03"
$ws.Range("D8").Value = "0.420"
$ws.Range("D9").Value = "7.06"
$ws.Range("D10").Value = "0.104"
$ws.Range("D11").Value = "0.347"
$ws.Range("D14").Value = "25.79"
$ws.Range("D15").Value = "0.0000155"
$ws.Range("D18").Value = "5.91"
$ws.Range("D19").Value = "12.35"
$ws.Range("D20").Value = "7.64"
$ws.Range("D21").Value = "313.76"
$ws.Range("D22").Value = "0.999"
$ws.Range("D23").Value = "5.75"
$ws.Range("D24").Value = "0.477"
$ws.Range("D25").Value = "62.17"
$ws.Range("D27").Value = "0.160"
$ws.Range("D29").Value = "6.35"
$ws.Range("D30").Value = "6.92"
$ws.Range("D31").Value = "1.73"
$ws.Range("D32").Value = "19.80"
$ws.Range("D33").Value = "1.11"
$ws.Range("D34").Value = "150.58"
$ws.Range("D35").Value = "4.40"
$ws.Range("D36").Value = "5.59"
$ws.Range("D38").Value = "23.47"
$ws.Range("D39").Value = "0.0644"
$ws.Range("D40").Value = "37.15"
$ws.Range("D43").Value = "3.64"
$ws.Range("D44").Value = "0.631"
$ws.Range("D46").Value = "1.32"
$ws.Range("D47").Value = "5.79"
$ws.Range("D48").Value = "0.895"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D50").Value = "18.66"
$ws.Range("D51").Value = "0.0840"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}

Write-Host "Cryptos list updated."
